$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Duplications" table (Table1) currently spans A13:F16. We need to
# add one more row (serial #5, given to "Ilya Ginzburg" / BADAS) so it
# grows to A13:F17.
$lo = $ws.ListObjects.Item("Table1")
$lo.ListRows.Add() | Out-Null

# Copy the formatting of the row above (row 16) down into the freshly
# added row 17 before we fill in values, so the new row keeps the same
# look (fonts/number formats/borders) as the rest of the table.
$ws.Range("A16:E16").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new duplication record.
$ws.Range("A17").Value = 5
$ws.Range("B17").Value = "Ilya Ginzburg"
$ws.Range("C17").Value = 45343
$ws.Range("D17").Value = "Ilya Ginzburg"
$ws.Range("E17").Value = 45343

# Reflect where the user ended up clicking next (just below the table).
$ws.Range("C18").Select()
